# Extend the register_1 bit-field table (rows 9-11) with two additional
# bit fields (bit_field_2 / bit_field_3), pushing every row below it down
# by two rows.
#
# Before (rows 8-11):
#   8  register_1 | 0x04 | bit_field_0 | 0:4  | ro
#   9             |      | bit_field_1 | 8:4  | ro
#  10  register_2 | 0x04 | bit_field_0 | 0:4  | wo | 0
#  11             |      | bit_field_1 | 8:4  | wo | 0
#
# After (rows 8-13):
#   8  register_1 | 0x04 | bit_field_0 | 0:4   | ro
#   9             |      | bit_field_1 | 8:4   | ro
#  10             |      | bit_field_2 | 16:8  | rof | 0xab
#  11             |      | bit_field_3 | 24:8  | reserved
#  12  register_2 | 0x04 | bit_field_0 | 0:4   | wo | 0
#  13             |      | bit_field_1 | 8:4   | wo | 0

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows right after row 9 (the current last row of the
# register_1 block), shifting every following row down by two.
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).Insert()

# Row 6 is a "middle" bit-field row (no bottom border) - reuse its
# formatting (borders/number formats) for the row-9/row-10 pair, which
# keeps the same "middle" look.
$ws.Range("B6:J6").Copy()
$ws.Range("B9:J10").PasteSpecial(-4122)

# Row 7 is the "last row of block" variant (adds the bottom border) -
# reuse that formatting for the new last row of the register_1 block.
$ws.Range("B7:J7").Copy()
$ws.Range("B11:J11").PasteSpecial(-4122)

# Row 9 keeps bit_field_1 (unchanged text, now on freshly styled cells).
$ws.Range("F9").Value = "bit_field_1"
$ws.Range("G9").Value = "8:4"
$ws.Range("H9").Value = "ro"

# Row 10: new bit_field_2, read-only-with-flag, initial value 0xab.
$ws.Range("F10").Value = "bit_field_2"
$ws.Range("G10").Value = "16:8"
$ws.Range("H10").Value = "rof"
$ws.Range("I10").Value = "0xab"

# Row 11: new bit_field_3, reserved.
$ws.Range("F11").Value = "bit_field_3"
$ws.Range("G11").Value = "24:8"
$ws.Range("H11").Value = "reserved"
